# fix: added guid to source system
#
# Adds a "ClassicSourceSystemGUID" text property to the ClassicSourceSystem
# view (Properties sheet) and registers the corresponding container
# ("ClassicSourceSystem", Used For = node) on the Containers sheet -
# mirroring the existing classicEquipmentGUID pattern used for
# ClassicEquipment.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Properties sheet: insert a new row right after the ClassicSourceSystem
#    "name" property row, describing the new ClassicSourceSystemGUID
#    property (mirrors the classicEquipmentGUID row further down).
# ---------------------------------------------------------------------
$wsProps = $wb.Worksheets.Item("Properties")
$wsProps.Rows.Item(54).Insert()

$wsProps.Cells.Item(54, 1).Value = "ClassicSourceSystem"        # View
$wsProps.Cells.Item(54, 2).Value = "ClassicSourceSystemGUID"    # View Property
$wsProps.Cells.Item(54, 6).Value = "text"                       # Value Type
$wsProps.Cells.Item(54, 7).Value = $true                        # Nullable
$wsProps.Cells.Item(54, 8).Value = $false                       # Immutable
$wsProps.Cells.Item(54, 9).Value = $false                       # Is List
$wsProps.Cells.Item(54, 11).Value = "ClassicSourceSystem"       # Container
$wsProps.Cells.Item(54, 12).Value = "ClassicSourceSystemGUID"   # Container Property

# ---------------------------------------------------------------------
# 2. Containers sheet: register the ClassicSourceSystem container
#    (same "Used For = node" pattern as the ClassicEquipment container).
# ---------------------------------------------------------------------
$wsContainers = $wb.Worksheets.Item("Containers")
$wsContainers.Cells.Item(4, 1).Value = "ClassicSourceSystem"    # Container
$wsContainers.Cells.Item(4, 5).Value = "node"                   # Used For

# ---------------------------------------------------------------------
# 3. Restore sensible selections and make Containers the active sheet,
#    matching where the author ended up after the edit.
# ---------------------------------------------------------------------
$wsProps.Activate() | Out-Null
$wsProps.Range("K54").Select() | Out-Null

$wsContainers.Activate() | Out-Null
$wsContainers.Range("F22").Select() | Out-Null
